$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new credentials row (username/password pair) after the existing data.
$ws.Range("A5").Value = "Admin"
$ws.Range("B5").Value = "admin123"

# Move the active selection to the newly written cell, matching Excel's
# behaviour of leaving the cursor on the last-edited cell.
[void]$ws.Range("B5").Select()
